# Updated cryptos list on Wed Jul 17 21:58:43 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.365.71"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.414.37"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'570.33"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "'156.46"
$ws.Range("E6").Value = "  -2.81%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = "  +8.35%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.419.32"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "'7.12"
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "4.001.16"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").Value = "'27.62"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "64.373.65"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "3.413.57"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "'13.85"
$ws.Range("E20").Value = "  -2.58%  "
$ws.Range("D21").Value = "'377.62"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "'8.00"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").Value = "'0.546"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'71.68"
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("D27").Value = "'10.32"
$ws.Range("E27").Value = "  +5.89%  "
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'1.47"
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("D34").Value = "'7.13"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("E35").Value = "  +7.26%  "
$ws.Range("D36").Value = "'160.09"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'6.97"
$ws.Range("E38").Value = "  +6.23%  "
$ws.Range("D39").Value = "'0.0762"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "2.879.61"
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("D41").Value = "'4.63"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").Value = "'26.32"
$ws.Range("E42").Value = "  -3.58%  "
$ws.Range("D43").Value = "'42.85"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'26.28"
$ws.Range("E44").Value = "  +5.75%  "
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "'320.43"
$ws.Range("E47").Value = "  +5.69%  "
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").Value = "'0.859"
$ws.Range("E51").Value = "  -2.41%  "
